# Update the three "Expected Utility" / "TRANSFER" shared-string values
# that moved for the waste-handling recompute (see commit "working on
# waste stuff").
#
#   B10 : "Expected Utility: 5.34296865663645"    -> "Expected Utility: 5.33770047422445"
#   B11 : "Expected Utility: 5.21822888677308"     -> "Expected Utility: 5.171968637050121"
#   L10 : "\n(TRANSFER self Foremz ((Electronics 1)) EU: -0.057523822012224654"
#            -> "\n(TRANSFER self Foremz ((Electronics 1)) EU: -0.43696381975762977"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B10").Value = "Expected Utility: 5.33770047422445"
$ws.Range("B11").Value = "Expected Utility: 5.171968637050121"
$ws.Range("L10").Value = "`n(TRANSFER self Foremz ((Electronics 1)) EU: -0.43696381975762977"
